# Patton's Best - Events.xlsx
# "Started to work on firing Sherman"
#
# Splits the old single "e053 Main Gun Firing" event into three events:
#   e053  - Main Gun Firing - Select Target   (rewritten text)
#   e053a - Main Gun Firing - No Target Available  (new)
#   e053b - Main Gun Firing (the old "select target / consult to-hit table" text, renamed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New text content (decoded from the authoritative diff)
# ---------------------------------------------------------------------------
$e053aShort = "e053a"
$e053bShort = "e053b"
$e053bText = "<Bold>e053b Main Gun Firing</Bold> `n<InlineUIContainer><Button Content='r4.74.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   `n<InlineUIContainer><Button Content='r9.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n<LineBreak/><LineBreak/>`nSelect a target by clicking the enemy unit on the Battle Board. `n<LineBreak/><LineBreak/>"
$e053SelectTargetText = "<Bold>e053 Main Gun Firing - Select Target</Bold> `n<InlineUIContainer><Button Content='r4.74.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   `n<InlineUIContainer><Button Content='r9.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n<LineBreak/><LineBreak/>`nSelect a target by clicking the enemy unit on the Battle Board. Only enemy units that have been spotted may be selected. The main gun may only fire at a target in the turret&apos;s sector unless the <Bold>Rotate Turret - Fire Main Gun</Bold> action was taken.`n<LineBreak/><LineBreak/>"
$e053aNoTargetText = "<Bold>e053a Main Gun Firing - No Target Available</Bold> `n<InlineUIContainer><Button Content='r4.74.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   `n<InlineUIContainer><Button Content='r9.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n<LineBreak/><LineBreak/>`nThere is no target that is selectable. Only enemy units that have been spotted may be selected. The main gun may only fire at a target in the turret&apos;s front sector unless the <Bold>Rotate Turret - Fire Main Gun</Bold> action was taken. Click image to continue.`n<LineBreak/><LineBreak/>`n                                            <InlineUIContainer><Image Name='Continue53a' Height='100' Width='100'></Image></InlineUIContainer>"

# ---------------------------------------------------------------------------
# 1) Row 63 (e053) gets new "Select Target" wording, and shrinks in height
# ---------------------------------------------------------------------------
$ws.Range("B63").Value2 = $e053SelectTargetText
$ws.Rows(63).RowHeight = 105

# ---------------------------------------------------------------------------
# 2) Insert two brand-new rows right after row 63 for e053a and e053b
# ---------------------------------------------------------------------------
$ws.Rows("64:65").Insert()

$ws.Range("A64").Value2 = $e053aShort
$ws.Range("B64").Value2 = $e053aNoTargetText
$ws.Rows(64).RowHeight = 120

$ws.Range("A65").Value2 = $e053bShort
$ws.Range("B65").Value2 = $e053bText
$ws.Rows(65).RowHeight = 90

# ---------------------------------------------------------------------------
# 3) Update the saved view state to match (best effort)
# ---------------------------------------------------------------------------
$ws.Range("B65").Select()
$excel.ActiveWindow.ScrollRow = 63
$excel.ActiveWindow.ScrollColumn = 1
